# Weekly price-sheet update: a new observation is inserted as row 63
# (dated 2023-06-05 / serial 45082), pushing the previously-existing
# rows 63..172 down to 64..173.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 63; this shifts rows 63:172 -> 64:173
# and extends the sheet dimension to A1:R173 automatically.
$ws.Rows("63:63").Insert()

# Populate the newly inserted row 63 with the new weekly data point.
$ws.Cells.Item(63, 1).Value = 10
$ws.Cells.Item(63, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(63, 3).Value = "La Araucanía"
$ws.Cells.Item(63, 4).Value = 45082
$ws.Cells.Item(63, 5).Value = 9
$ws.Cells.Item(63, 6).Value = 100112035
$ws.Cells.Item(63, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 80
$ws.Cells.Item(63, 11).Value = 28000
$ws.Cells.Item(63, 12).Value = 28000
$ws.Cells.Item(63, 13).Value = 28000
$ws.Cells.Item(63, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(63, 15).Value = "Región Metropolitana"
$ws.Cells.Item(63, 16).Value = 1867
$ws.Cells.Item(63, 17).Value = 15
$ws.Cells.Item(63, 18).Value = "Hortaliza"
